$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: Roi, 14.1, 5h
$ws.Range("B7").Value = "Roi"
$ws.Range("C7").Value = 14.1
$ws.Range("D7").Value = "5h"

# Row 8: Roi, 15.1, 8h
$ws.Range("B8").Value = "Roi"
$ws.Range("C8").Value = 15.1
$ws.Range("D8").Value = "8h"

# Row 9: Magal, 15.1, 8h
$ws.Range("B9").Value = "Magal"
$ws.Range("C9").Value = 15.1
$ws.Range("D9").Value = "8h"

# Row 10: Shay, 15.1, 8h
$ws.Range("B10").Value = "Shay"
$ws.Range("C10").Value = 15.1
$ws.Range("D10").Value = "8h"

# Row 11: Roi, 16.1, 5h
$ws.Range("B11").Value = "Roi"
$ws.Range("C11").Value = 16.1
$ws.Range("D11").Value = "5h"

$ws.Range("D11").Select()
